$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "SamplePojo"

# Pre-format the date (text) columns so the auto date-detection doesn't
# convert the "dd.mm.yyyy" looking strings into real date serials.
$dateCells = $ws.Range("G1:G5")
$dateCells.NumberFormat = "@"

# Row 1
$ws.Range("A1").Value = 1
$ws.Range("B1").Value = "Seray"
$ws.Range("C1").Value = "Uzgur"
$ws.Range("D1").Value = 11
$ws.Range("E1").Value = 111
$ws.Range("F1").Value = 1111
$ws.Range("G1").Value = "01.01.2014"

# Row 2
$ws.Range("A2").Value = 2
$ws.Range("B2").Value = "Kaan"
$ws.Range("C2").Value = "Alkim"
$ws.Range("D2").Value = 12
$ws.Range("E2").Value = 112
$ws.Range("F2").Value = 1112
$ws.Range("G2").Value = "02.01.2014"

# Row 3
$ws.Range("A3").Value = 3
$ws.Range("B3").Value = "Sinan"
$ws.Range("C3").Value = "Selimogli"
$ws.Range("D3").Value = 13
$ws.Range("E3").Value = 113
$ws.Range("F3").Value = 1113
$ws.Range("G3").Value = "03.01.2014"

# Row 4
$ws.Range("A4").Value = 4
$ws.Range("B4").Value = "Kamil"
$ws.Range("C4").Value = "Bukum"
$ws.Range("D4").Value = 14
$ws.Range("E4").Value = 114
$ws.Range("F4").Value = 1114
$ws.Range("G4").Value = "04.01.2014"

# Row 5
$ws.Range("A5").Value = 5
$ws.Range("B5").Value = "Hasan"
$ws.Range("C5").Value = "Mumin"
$ws.Range("D5").Value = 15
$ws.Range("E5").Value = 115
$ws.Range("F5").Value = 1115
$ws.Range("G5").Value = "05.01.2014"

# Remove the temporary text format again so the cells keep the
# workbook's default (General) style.
$dateCells.ClearFormats()

# Column widths (A, B, G) as in the target layout.
$ws.Columns.Item(1).ColumnWidth = 12 - 0.8333333333333339
$ws.Columns.Item(2).ColumnWidth = 10 - 0.8333333333333339
$ws.Columns.Item(7).ColumnWidth = 9.6640625 - 0.8333333333333339

# Select the header row, matching the saved selection in the workbook.
$ws.Range("A1:XFD1").Select()

# Page setup tweaks.
$ws.PageSetup.Orientation = 1
